## Applies the bill_template.docx edit described by the commit:
## "Implement bill generation endpoint with dynamic item totals and discount calculations"
##
## Concretely, inside the single invoice line-items table:
##   - header row: the 3-run {price} placeholder is merged into one run
##   - the "#items" row: "#items" becomes "/items" (the leading '#' is swapped for '/')
##   - the "#items" row: {qty1} and {price1} become empty {}
##   - the "#items" row: the literal "100" amount cell is cleared (text + center
##     alignment removed, leaving an empty paragraph)
##   - the four "item N" rows (item 4..7): {item N}, {qty1}, {price1} all become
##     empty {} (their "100" amount cells are left untouched)

$d = $word.ActiveDocument

function Find-NthRange {
    # Returns a Range positioned on the n-th (1-based) occurrence of $searchText
    # in document order, without mutating anything.
    param(
        [object]$doc,
        [string]$searchText,
        [int]$n
    )
    $found = $doc.Content
    for ($i = 1; $i -le $n; $i++) {
        $ok = $found.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
        if (-not $ok) { return $null }
        if ($i -lt $n) {
            $found = $doc.Range($found.End, $doc.Content.End)
        }
    }
    return $found
}

function Replace-NthOccurrence {
    # Replaces only the n-th (1-based) occurrence of $searchText in the whole
    # document with $replacement (single shot - wdReplaceOne), leaving every
    # other occurrence (earlier or later) untouched.
    param(
        [object]$doc,
        [string]$searchText,
        [int]$n,
        [string]$replacement
    )
    $target = Find-NthRange $doc $searchText $n
    if ($null -eq $target) {
        throw "Could not find occurrence #$n of '$searchText'"
    }
    $target.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, $replacement, 1) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) Header row: "{" + "price" + "}" (3 runs) -> single run "{price}"
#    "{price}" is unique in the document, so a plain replace-all is safe and
#    Word naturally coalesces the matched runs into one.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("{price}", $true, $false, $false, $false, $false, $true, 0, $false, "{price}", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "#items" literal text. The exact run-text "#items" occurs twice in the
#    document: once in an explanatory legend row (must stay untouched) and
#    once in the actual data row we need to edit (the 2nd occurrence).
# ---------------------------------------------------------------------------
Replace-NthOccurrence $d "#items" 2 "/items"

# ---------------------------------------------------------------------------
# 3) The five "{qty1}" cells: row "#items", then rows item4..item7.
#    Only the first one (row "#items") is touched by this commit.
# ---------------------------------------------------------------------------
Replace-NthOccurrence $d "{qty1}" 1 "{}"

# ---------------------------------------------------------------------------
# 4) The five "{price1}" cells - again only the first (row "#items").
# ---------------------------------------------------------------------------
Replace-NthOccurrence $d "{price1}" 1 "{}"

# ---------------------------------------------------------------------------
# 5) The "100" amount cell belonging to the same "#items" row: remove the
#    center alignment and delete the run/text entirely (empty paragraph).
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)
$itemsRow = $table.Rows.Item(10)
$amountCell = $itemsRow.Cells.Item(6)
$amountCell.Range.Paragraphs.Item(1).Alignment = 0
$amountCell.Range.Find.Execute("100", $true, $false, $false, $false, $false, $true, 0, $false, "", 1) | Out-Null

# ---------------------------------------------------------------------------
# 6) The four "item N" rows (N = 4,5,6,7): "{item N}" -> "{}" and their
#    {qty1}/{price1} -> "{}" too (their "100" cell is left as-is).
# ---------------------------------------------------------------------------
foreach ($n in 4, 5, 6, 7) {
    $placeholder = "{item $n}"
    $d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 0, $false, "{}", 2) | Out-Null

    # after every previous step, the remaining first occurrence of {qty1} /
    # {price1} is always the one belonging to this row
    Replace-NthOccurrence $d "{qty1}" 1 "{}"
    Replace-NthOccurrence $d "{price1}" 1 "{}"
}

Write-Output "done"
